# Rename the "index" column (table header / shared string) to "i".
# Setting the header cell value updates the shared-strings table, the
# worksheet cell, and the associated table column name together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "i"

# Convert the index column from 1-based to 0-based values (decrement
# every data row in column A by 1).
for ($r = 2; $r -le 503; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 - 1
}

# Narrow column A now that it holds shorter values / a shorter header.
$ws.Columns.Item(1).ColumnWidth = 3.17
